# Jellyfish stats workbook update
# - Adds target_bitrate_kbps / achieved_bitrate_kbps columns (B, C)
# - Renames / shifts the prediction / residual / other columns to
#   prediction_kbps / residual_kbps / other_kbps (D, E, F) with new values
# - Adds avg_psnr / avg_qp columns (G, H)
# - Removes the old duplicate F:I block
# - Repoints both charts at the new columns/categories

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row
# ---------------------------------------------------------------------------
$headers = New-Object 'object[,]' 1,8
$headers[0,0] = "total"
$headers[0,1] = "target_bitrate_kbps"
$headers[0,2] = "achieved_bitrate_kbps"
$headers[0,3] = "prediction_kbps"
$headers[0,4] = "residual_kbps"
$headers[0,5] = "other_kbps"
$headers[0,6] = "avg_psnr"
$headers[0,7] = "avg_qp"
$ws.Range("A1:H1").Value = $headers

# ---------------------------------------------------------------------------
# 2. Data rows 2-11 (A untouched values repeated here for completeness)
# ---------------------------------------------------------------------------
$arr = New-Object 'object[,]' 10,8

$arr[0,0] = 363976743
$arr[0,1] = 100000
$arr[0,2] = 96963.404335200001
$arr[0,3] = 5768.3638620000002
$arr[0,4] = 90730.9230885
$arr[0,5] = 464.1173847
$arr[0,6] = 32.050589000000002
$arr[0,7] = 6.29

$arr[1,0] = 289970358
$arr[1,1] = 80000
$arr[1,2] = 77248.103371200006
$arr[1,3] = 5824.6339688999997
$arr[1,4] = 70931.601462599996
$arr[1,5] = 491.867939699999
$arr[1,6] = 32.047114999999998
$arr[1,7] = 8.14

$arr[2,0] = 217022706
$arr[2,1] = 60000
$arr[2,2] = 57814.8488784
$arr[2,3] = 5507.6819381999903
$arr[2,4] = 51789.377148300002
$arr[2,5] = 517.78979189999995
$arr[2,6] = 32.045279999999998
$arr[2,7] = 10.32

$arr[3,0] = 144913374
$arr[3,1] = 40000
$arr[3,2] = 38604.922833600001
$arr[3,3] = 4769.0103491999998
$arr[3,4] = 33317.419130099901
$arr[3,5] = 518.49335429999996
$arr[3,6] = 32.040278000000001
$arr[3,7] = 13.41

$arr[4,0] = 73105843
$arr[4,1] = 20000
$arr[4,2] = 19475.396575199899
$arr[4,3] = 3462.0723287999999
$arr[4,4] = 15558.3974619
$arr[4,5] = 454.9267845
$arr[4,6] = 32.014301000000003
$arr[4,7] = 18.79

$arr[5,0] = 36544663
$arr[5,1] = 10000
$arr[5,2] = 9735.4982232000002
$arr[5,3] = 2352.8403045
$arr[5,4] = 7034.9366547
$arr[5,5] = 347.72126400000002
$arr[5,6] = 31.947144999999999
$arr[5,7] = 24.47

$arr[6,0] = 18263199
$arr[6,1] = 5000
$arr[6,2] = 4865.3162136000001
$arr[6,3] = 1476.6396489000001
$arr[6,4] = 3147.1835993999998
$arr[6,5] = 241.49296530000001
$arr[6,6] = 31.799208
$arr[6,7] = 30.41

$arr[7,0] = 3707878
$arr[7,1] = 1000
$arr[7,2] = 987.77869920000001
$arr[7,3] = 395.27692739999998
$arr[7,4] = 492.7139262
$arr[7,5] = 99.787845599999997
$arr[7,6] = 30.975273999999999
$arr[7,7] = 42.2

$arr[8,0] = 1882402
$arr[8,1] = 500
$arr[8,2] = 501.47189279999998
$arr[8,3] = 216.0970533
$arr[8,4] = 211.97571209999899
$arr[8,5] = 73.399127399999998
$arr[8,6] = 30.377583000000001
$arr[8,7] = 46.64

$arr[9,0] = 1106036
$arr[9,1] = 250
$arr[9,2] = 294.64799040000003
$arr[9,3] = 160.797708
$arr[9,4] = 101.7010305
$arr[9,5] = 32.149251900000003
$arr[9,6] = 28.673165999999998
$arr[9,7] = 51

$ws.Range("A2:H11").Value = $arr

# ---------------------------------------------------------------------------
# 3. Remove the old leftover duplicate block (F:I) beyond the new H column
# ---------------------------------------------------------------------------
$ws.Range("I1:I11").Clear()

# ---------------------------------------------------------------------------
# 4. Repoint chart 1 (full range, rows 2-11) series at the new columns
# ---------------------------------------------------------------------------
$chart1 = $ws.ChartObjects().Item(1).Chart

$s1 = $chart1.SeriesCollection().Item(1)
$s1.Name = "=Sheet1!`$D`$1"
$s1.Values = $ws.Range("D2:D11")
$s1.XValues = $ws.Range("B2:B11")

$s2 = $chart1.SeriesCollection().Item(2)
$s2.Name = "=Sheet1!`$E`$1"
$s2.Values = $ws.Range("E2:E11")
$s2.XValues = $ws.Range("B2:B11")

$s3 = $chart1.SeriesCollection().Item(3)
$s3.Name = "=Sheet1!`$F`$1"
$s3.Values = $ws.Range("F2:F11")
$s3.XValues = $ws.Range("B2:B11")

$chart1.ChartGroups(1).GapWidth = 150

# ---------------------------------------------------------------------------
# 5. Repoint chart 2 (zoomed range, rows 8-11) series at the new columns
# ---------------------------------------------------------------------------
$chart2 = $ws.ChartObjects().Item(2).Chart

$t1 = $chart2.SeriesCollection().Item(1)
$t1.Name = "=Sheet1!`$D`$1"
$t1.Values = $ws.Range("D8:D11")
$t1.XValues = $ws.Range("B8:B11")

$t2 = $chart2.SeriesCollection().Item(2)
$t2.Name = "=Sheet1!`$E`$1"
$t2.Values = $ws.Range("E8:E11")
$t2.XValues = $ws.Range("B8:B11")

$t3 = $chart2.SeriesCollection().Item(3)
$t3.Name = "=Sheet1!`$F`$1"
$t3.Values = $ws.Range("F8:F11")
$t3.XValues = $ws.Range("B8:B11")

$wb.RefreshAll()
$excel.Calculate()
